$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L (12), shifting old L..P to M..Q
$ws.Columns.Item(12).Insert()

# New column L (12) header/value -- new header cell keeps the generic header style (s=1,
# same as the other header cells), new data cell inherits the left neighbour's (K2) style (s=5)
$ws.Cells.Item(1, 12).Value = "SIDEBAR_SUBMENU_SUBMENU"
$ws.Cells.Item(2, 12).Value = "Setup Kelengkapan Kepesertaan"

$ws.Cells.Item(2, 11).Copy()
$ws.Cells.Item(2, 12).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set new column width to match column K's width (15)
$ws.Columns.Item(12).ColumnWidth = $ws.Columns.Item(11).ColumnWidth

# Update selection
$ws.Range("O12").Select()
